# All_F1_tables.xlsx edit:
#  - Sheet3 (100 Australian species / COI): the "DADA2Tax" row is removed and
#    every remaining data row gets fresh Precision/Recall/F1/F0.5/Accuracy
#    numbers (full rewrite of the COI part).
#  - Sheet6 (was Lutjanidae/COI) is overwritten with what used to be Sheet7's
#    data (Wadjemup/12S).
#  - Sheet7 (was Wadjemup/12S) is overwritten with what used to be Sheet8's
#    data (Wadjemup/16S).
#  - Sheet8 (Wadjemup/16S) and Sheet9 (Wadjemup/COI) are deleted outright.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Copy Sheet7's (still-original Wadjemup/12S) data into Sheet6 first, then
# copy Sheet8's (Wadjemup/16S) data into Sheet7.
for ($r = 2; $r -le 15; $r++) {
    for ($c = 1; $c -le 8; $c++) {
        $wb.Worksheets.Item("Sheet6").Cells.Item($r, $c).Value2 = $wb.Worksheets.Item("Sheet7").Cells.Item($r, $c).Value2
    }
}

for ($r = 2; $r -le 15; $r++) {
    for ($c = 1; $c -le 8; $c++) {
        $wb.Worksheets.Item("Sheet7").Cells.Item($r, $c).Value2 = $wb.Worksheets.Item("Sheet8").Cells.Item($r, $c).Value2
    }
}

# Drop the now-redundant sheets (highest index first -- fresh lookups, no
# cached worksheet references, to avoid stale COM handles after a delete).
$wb.Worksheets.Item("Sheet9").Delete() | Out-Null
$wb.Worksheets.Item("Sheet8").Delete() | Out-Null

# Rewrite Sheet3 (100 Australian species / COI): remove the "DADA2Tax" row
# (currently row 6) which shifts everything below it up by one, then write
# the refreshed metrics for every method.
$sheet3 = $wb.Worksheets.Item("Sheet3")
$sheet3.Rows.Item(6).Delete() | Out-Null

$sheet3Data = @(
    @("Kraken2_0.0", 0, 0, 0, 0, 0.1515151515151515),
    @("Kraken2_0.05", 0, 0, 0, 0, 0.1515151515151515),
    @("Kraken2_0.1", 0, 0, 0, 0, 0.1515151515151515)
)

$sheet3.Cells.Item(2, 4).Value2 = 1
$sheet3.Cells.Item(2, 5).Value2 = 0.6904761904761905
$sheet3.Cells.Item(2, 6).Value2 = 0.8169014084507042
$sheet3.Cells.Item(2, 7).Value2 = 0.9177215189873418
$sheet3.Cells.Item(2, 8).Value2 = 0.7373737373737373

$sheet3.Cells.Item(3, 4).Value2 = 0.9803921568627451
$sheet3.Cells.Item(3, 5).Value2 = 0.6024096385542169
$sheet3.Cells.Item(3, 6).Value2 = 0.746268656716418
$sheet3.Cells.Item(3, 7).Value2 = 0.8710801393728221
$sheet3.Cells.Item(3, 8).Value2 = 0.6565656565656566

$sheet3.Cells.Item(4, 4).Value2 = 0.863013698630137
$sheet3.Cells.Item(4, 5).Value2 = 0.8513513513513513
$sheet3.Cells.Item(4, 6).Value2 = 0.8571428571428572
$sheet3.Cells.Item(4, 7).Value2 = 0.8606557377049179
$sheet3.Cells.Item(4, 8).Value2 = 0.7878787878787878

$sheet3.Cells.Item(5, 4).Value2 = 1
$sheet3.Cells.Item(5, 5).Value2 = 0.6785714285714286
$sheet3.Cells.Item(5, 6).Value2 = 0.8085106382978724
$sheet3.Cells.Item(5, 7).Value2 = 0.9134615384615385
$sheet3.Cells.Item(5, 8).Value2 = 0.7272727272727273

# Row 6 is now "Kraken2_0.0" (DADA2Tax removed), row 7 "Kraken2_0.05", row 8
# "Kraken2_0.1" -- labels unchanged, only the metrics change.
$sheet3.Cells.Item(6, 4).Value2 = 0
$sheet3.Cells.Item(6, 5).Value2 = 0
$sheet3.Cells.Item(6, 6).Value2 = 0
$sheet3.Cells.Item(6, 7).Value2 = 0
$sheet3.Cells.Item(6, 8).Value2 = 0.1515151515151515

$sheet3.Cells.Item(7, 4).Value2 = 0
$sheet3.Cells.Item(7, 5).Value2 = 0
$sheet3.Cells.Item(7, 6).Value2 = 0
$sheet3.Cells.Item(7, 7).Value2 = 0
$sheet3.Cells.Item(7, 8).Value2 = 0.1515151515151515

$sheet3.Cells.Item(8, 4).Value2 = 0
$sheet3.Cells.Item(8, 5).Value2 = 0
$sheet3.Cells.Item(8, 6).Value2 = 0
$sheet3.Cells.Item(8, 7).Value2 = 0
$sheet3.Cells.Item(8, 8).Value2 = 0.1515151515151515

$sheet3.Cells.Item(9, 4).Value2 = 1
$sheet3.Cells.Item(9, 5).Value2 = 0.6785714285714286
$sheet3.Cells.Item(9, 6).Value2 = 0.8085106382978724
$sheet3.Cells.Item(9, 7).Value2 = 0.9134615384615385
$sheet3.Cells.Item(9, 8).Value2 = 0.7272727272727273

$sheet3.Cells.Item(10, 4).Value2 = 0.9682539682539683
$sheet3.Cells.Item(10, 5).Value2 = 0.7439024390243902
$sheet3.Cells.Item(10, 6).Value2 = 0.8413793103448277
$sheet3.Cells.Item(10, 7).Value2 = 0.9131736526946107
$sheet3.Cells.Item(10, 8).Value2 = 0.7676767676767676

$sheet3.Cells.Item(11, 4).Value2 = 0.9230769230769231
$sheet3.Cells.Item(11, 5).Value2 = 0.2926829268292683
$sheet3.Cells.Item(11, 6).Value2 = 0.4444444444444444
$sheet3.Cells.Item(11, 7).Value2 = 0.6451612903225807
$sheet3.Cells.Item(11, 8).Value2 = 0.3939393939393939

$sheet3.Cells.Item(12, 4).Value2 = 0.8648648648648649
$sheet3.Cells.Item(12, 5).Value2 = 0.8648648648648649
$sheet3.Cells.Item(12, 6).Value2 = 0.8648648648648649
$sheet3.Cells.Item(12, 7).Value2 = 0.8648648648648649
$sheet3.Cells.Item(12, 8).Value2 = 0.797979797979798

$sheet3.Cells.Item(13, 4).Value2 = 0.864406779661017
$sheet3.Cells.Item(13, 5).Value2 = 0.6219512195121951
$sheet3.Cells.Item(13, 6).Value2 = 0.723404255319149
$sheet3.Cells.Item(13, 7).Value2 = 0.8018867924528302
$sheet3.Cells.Item(13, 8).Value2 = 0.6060606060606061

$sheet3.Cells.Item(14, 4).Value2 = 0.8695652173913043
$sheet3.Cells.Item(14, 5).Value2 = 0.8
$sheet3.Cells.Item(14, 6).Value2 = 0.8333333333333333
$sheet3.Cells.Item(14, 7).Value2 = 0.8547008547008546
$sheet3.Cells.Item(14, 8).Value2 = 0.7575757575757576

Write-Host "edit complete"
